$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1) — replace the old TBD-participant columns with the new
# "Mixed ANOVA" data-collection header set.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Participant"
$ws.Range("B1").Value = "Order"
$ws.Range("C1").Value = "Course"
$ws.Range("D1").Value = "Time"
$ws.Range("E1").Value = "Kills"
$ws.Range("F1").Value = "Shots_Fired"
$ws.Range("G1").Value = "Shots_Hit"
$ws.Range("H1").Value = "Hit_Percent"

# ---------------------------------------------------------------------------
# Data rows 2-9.
# Columns:
#   A Participant #   B Order (A->B / B->A)   C Course (A / B)
#   D Time   E Kills   F Shots_Fired   G Shots_Hit   H Hit_Percent
# ---------------------------------------------------------------------------

# Row 2
$ws.Range("A2").Formula = "=1"
$ws.Range("B2").Value = "A->B"
$ws.Range("C2").Value = "A"
$ws.Range("D2").Value = 10.5
$ws.Range("E2").Formula = "=1"
$ws.Range("F2").Value = 10
$ws.Range("G2").Formula = "=1"
$ws.Range("H2").Formula = "=G2/F2*100"

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "A->B"
$ws.Range("C3").Value = "B"
$ws.Range("D3").Formula = "=D2+1"
$ws.Range("E3").Formula = "=E2+1"
$ws.Range("F3").Value = 15
$ws.Range("G3").Formula = "=G2+1"
$ws.Range("H3").Formula = "=G3/F3*100"

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "B->A"
$ws.Range("C4").Value = "B"
$ws.Range("D4").Formula = "=D3+1"
$ws.Range("E4").Formula = "=E3+1"
$ws.Range("F4").Formula = "=19"
$ws.Range("G4").Formula = "=G3+1"
$ws.Range("H4").Formula = "=G4/F4*100"

# Row 5
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "B->A"
$ws.Range("C5").Value = "A"
$ws.Range("D5").Formula = "=D4+1"
$ws.Range("E5").Formula = "=E4+1"
$ws.Range("F5").Formula = "=19"
$ws.Range("G5").Formula = "=G4+1"
$ws.Range("H5").Formula = "=G5/F5*100"

# Row 6
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "A->B"
$ws.Range("C6").Value = "A"
$ws.Range("D6").Formula = "=D5+1"
$ws.Range("E6").Formula = "=E5+1"
$ws.Range("F6").Formula = "=19"
$ws.Range("G6").Formula = "=1"
$ws.Range("H6").Formula = "=G6/F6*100"

# Row 7
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "A->B"
$ws.Range("C7").Value = "B"
$ws.Range("D7").Formula = "=D6+1"
$ws.Range("E7").Formula = "=E6+1"
$ws.Range("F7").Formula = "=19"
$ws.Range("G7").Formula = "=G6+1"
$ws.Range("H7").Formula = "=G7/F7*100"

# Row 8
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "B->A"
$ws.Range("C8").Value = "B"
$ws.Range("D8").Formula = "=D7+1"
$ws.Range("E8").Formula = "=E7+1"
$ws.Range("F8").Formula = "=19"
$ws.Range("G8").Formula = "=G7+1"
$ws.Range("H8").Formula = "=G8/F8*100"

# Row 9
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "B->A"
$ws.Range("C9").Value = "A"
$ws.Range("D9").Formula = "=D8+1"
$ws.Range("E9").Formula = "=E8+1"
$ws.Range("F9").Formula = "=19"
$ws.Range("G9").Formula = "=G8+1"
$ws.Range("H9").Formula = "=G9/F9*100"

# ---------------------------------------------------------------------------
# Column widths for the newly-added columns F, G, H (best-fit / custom).
# (ColumnWidth values chosen so the engine's pixel-quantised stored width
# lands on the same bucket as the authored widths of 10.44140625 / 9.109375
# / 10.88671875 characters.)
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 9.666666666666666
$ws.Columns.Item(7).ColumnWidth = 8.333333333333334
$ws.Columns.Item(8).ColumnWidth = 10.0

# Move the active-cell selection like the author's session ended up.
$ws.Range("N14").Select()

Write-Output "done"
